$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCells = @("D2","D3","D5","D7","D8","D11","D12","D13","D14","D15","D16","D18","D19","D22","D25","D26","D31","D35","D36","D41","D42","D43","D45","D47","D49","D50")
foreach ($ref in $dCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = '33.899.61'
$ws.Range("E2").Value = '  -0.84%  '
$ws.Range("D3").Value = '1.781.07'
$ws.Range("E3").Value = '  -1.31%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '224.31'
$ws.Range("E5").Value = '  +0.45%  '
$ws.Range("E6").Value = '  -1.33%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").Value = '31.80'
$ws.Range("E8").Value = '  -3.13%  '
$ws.Range("E9").Value = '  +0.54%  '
$ws.Range("E10").Value = '  -5.63%  '
$ws.Range("D11").Value = '0.0937'
$ws.Range("E11").Value = '  +0.84%  '
$ws.Range("D12").Value = '2.037.28'
$ws.Range("E12").Value = '  -1.31%  '
$ws.Range("D13").Value = '1.837.34'
$ws.Range("E13").Value = '  +1.77%  '
$ws.Range("D14").Value = '11.16'
$ws.Range("E14").Value = '  +1.73%  '
$ws.Range("D15").Value = '33.882.47'
$ws.Range("E15").Value = '  -0.97%  '
$ws.Range("D16").Value = '0.610'
$ws.Range("E16").Value = '  -3.58%  '
$ws.Range("E17").Value = '  -2.28%  '
$ws.Range("D18").Value = '66.99'
$ws.Range("E18").Value = '  -2.84%  '
$ws.Range("D19").Value = '239.43'
$ws.Range("E19").Value = '  -3.28%  '
$ws.Range("E20").Value = '  -2.39%  '
$ws.Range("E21").Value = '  +0.02%  '
$ws.Range("D22").Value = '10.57'
$ws.Range("E22").Value = '  -5.27%  '
$ws.Range("E23").Value = '  -1.94%  '
$ws.Range("E24").Value = '  -3.11%  '
$ws.Range("D25").Value = '161.00'
$ws.Range("E25").Value = '  +0.77%  '
$ws.Range("D26").Value = '16.11'
$ws.Range("E26").Value = '  -2.97%  '
$ws.Range("E27").Value = '  -1.92%  '
$ws.Range("E28").Value = '  -0.95%  '
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("E30").Value = '  +0.89%  '
$ws.Range("D31").Value = '0.0510'
$ws.Range("E31").Value = '  -3.95%  '
$ws.Range("E32").Value = '  -3.92%  '
$ws.Range("E33").Value = '  -0.61%  '
$ws.Range("E34").Value = '  -1.94%  '
$ws.Range("D35").Value = '1.392.51'
$ws.Range("E35").Value = '  -1.82%  '
$ws.Range("D36").Value = '0.637'
$ws.Range("E36").Value = '  -2.84%  '
$ws.Range("E37").Value = '  -1.57%  '
$ws.Range("E38").Value = '  -1.09%  '
$ws.Range("E39").Value = '  +4.54%  '
$ws.Range("E40").Value = '  -0.03%  '
$ws.Range("D41").Value = '0.916'
$ws.Range("E41").Value = '  -2.37%  '
$ws.Range("D42").Value = '78.45'
$ws.Range("E42").Value = '  -3.00%  '
$ws.Range("D43").Value = '13.55'
$ws.Range("E43").Value = '  +12.39%  '
$ws.Range("E44").Value = '  -3.12%  '
$ws.Range("D45").Value = '0.0508'
$ws.Range("E45").Value = '  +2.19%  '
$ws.Range("E46").Value = '  +1.85%  '
$ws.Range("D47").Value = '0.0₆0134'
$ws.Range("E47").Value = '  +6.78%  '
$ws.Range("E48").Value = '  -1.77%  '
$ws.Range("D49").Value = '1.938.77'
$ws.Range("E49").Value = '  -1.22%  '
$ws.Range("D50").Value = '105.45'
$ws.Range("E50").Value = '  -2.84%  '
$ws.Range("E51").Value = '  -0.06%  '
